$wb = $excel.ActiveWorkbook

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1667.4166
$ws.Range("I80").Value = 2573
$ws.Range("J80").Value = 399.6
$ws.Range("K80").Value = 7719
$ws.Range("L80").Value = 1198.8
$ws.Range("M80").Value = -6721
$ws.Range("N80").Value = -3194.8

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1667.4166
$ws.Range("I83").Value = 2573
$ws.Range("J83").Value = 399.6
$ws.Range("K83").Value = 23157
$ws.Range("L83").Value = 3596.4
$ws.Range("M83").Value = -18165
$ws.Range("N83").Value = -13580.4

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 535.2
$ws.Range("I107").Value = 554.6667
$ws.Range("K107").Value = 554.6667
$ws.Range("M107").Value = 1365.3333

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6945.9697
$ws.Range("I116").Value = 6909.607
$ws.Range("J116").Value = 7149.6
$ws.Range("K116").Value = 6909.607
$ws.Range("L116").Value = 7149.6
$ws.Range("M116").Value = -3467.607
$ws.Range("N116").Value = -14033.6

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1469.75
$ws.Range("I118").Value = 1469.75
$ws.Range("K118").Value = 4409.25
$ws.Range("M118").Value = -2752.25

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4328.3076
$ws.Range("I132").Value = 4372.3335
$ws.Range("K132").Value = 13117.0005
$ws.Range("M132").Value = -10587.0005

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 45454960
$ws.Range("I135").Value = 55555988
$ws.Range("K135").Value = 500003892
$ws.Range("M135").Value = -500001357

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2578.25
$ws.Range("I141").Value = 1771.3334
$ws.Range("K141").Value = 5314.0002
$ws.Range("M141").Value = -134.0002000000004

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2322.2104
$ws.Range("I2").Value = 1968.8334
$ws.Range("K2").Value = 1968.8334
$ws.Range("M2").Value = -1855.8334

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3311.2292
$ws.Range("I32").Value = 2150.8696
$ws.Range("K32").Value = 2150.8696
$ws.Range("M32").Value = -1863.8696

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5109.3887
$ws.Range("I45").Value = 8445.25
$ws.Range("K45").Value = 8445.25
$ws.Range("M45").Value = -8068.25

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 29416808
$ws.Range("I74").Value = 37042356
$ws.Range("J74").Value = 3982
$ws.Range("K74").Value = 37042356
$ws.Range("L74").Value = 3982
$ws.Range("M74").Value = -37041482
$ws.Range("N74").Value = -5730

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 29416808
$ws.Range("I77").Value = 37042356
$ws.Range("J77").Value = 3982
$ws.Range("K77").Value = 185211780
$ws.Range("L77").Value = 19910
$ws.Range("M77").Value = -185207412
$ws.Range("N77").Value = -28646

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2322.2104
$ws.Range("I116").Value = 1968.8334
$ws.Range("K116").Value = 1968.8334
$ws.Range("M116").Value = 325.1666

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5633.04
$ws.Range("I122").Value = 5022.684
$ws.Range("J122").Value = 7565.8335
$ws.Range("K122").Value = 15068.052
$ws.Range("L122").Value = 22697.5005
$ws.Range("M122").Value = -12618.052
$ws.Range("N122").Value = -27597.5005

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2322.2104
$ws.Range("I3").Value = 1968.8334
$ws.Range("K3").Value = 1968.8334
$ws.Range("M3").Value = -1854.8334

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 71005.13
$ws.Range("I107").Value = 4314.8335
$ws.Range("J107").Value = 337766.34
$ws.Range("K107").Value = 4314.8335
$ws.Range("L107").Value = 337766.34
$ws.Range("M107").Value = -2394.8335
$ws.Range("N107").Value = -341606.34

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17244750
$ws.Range("I134").Value = 20003378
$ws.Range("K134").Value = 60010134
$ws.Range("M134").Value = -60007599

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11617.904
$ws.Range("I86").Value = 8986.666999999999
$ws.Range("J86").Value = 13591.333
$ws.Range("K86").Value = 8986.666999999999
$ws.Range("L86").Value = 13591.333
$ws.Range("M86").Value = -7863.666999999999
$ws.Range("N86").Value = -15837.333

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11617.904
$ws.Range("I89").Value = 8986.666999999999
$ws.Range("J89").Value = 13591.333
$ws.Range("K89").Value = 44933.335
$ws.Range("L89").Value = 67956.66500000001
$ws.Range("M89").Value = -39317.335
$ws.Range("N89").Value = -79188.66500000001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9302290
$ws.Range("I134").Value = 10045935
$ws.Range("K134").Value = 30137805
$ws.Range("M134").Value = -30135270

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 546.1875
$ws.Range("I33").Value = 325.77777
$ws.Range("J33").Value = 829.5714
$ws.Range("K33").Value = 1954.66662
$ws.Range("L33").Value = 4977.428400000001
$ws.Range("M33").Value = -1671.66662
$ws.Range("N33").Value = -5543.428400000001

# CUL row 123
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 199.5
$ws.Range("I123").Value = 199.5
$ws.Range("K123").Value = 598.5
$ws.Range("M123").Value = 1851.5

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 9999
$ws.Range("J125").Value = 9999
$ws.Range("L125").Value = 29997
$ws.Range("N125").Value = -39837

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1790.6428
$ws.Range("I131").Value = 1168.5
$ws.Range("J131").Value = 2039.5
$ws.Range("K131").Value = 3505.5
$ws.Range("L131").Value = 6118.5
$ws.Range("M131").Value = 1534.5
$ws.Range("N131").Value = -16198.5

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 185.8
$ws.Range("I2").Value = 162.5
$ws.Range("J2").Value = 201.33333
$ws.Range("K2").Value = 162.5
$ws.Range("L2").Value = 201.33333
$ws.Range("M2").Value = -49.5
$ws.Range("N2").Value = -427.33333

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3166.3572
$ws.Range("I102").Value = 3435
$ws.Range("J102").Value = 1554.5
$ws.Range("K102").Value = 3435
$ws.Range("L102").Value = 1554.5
$ws.Range("M102").Value = -1813
$ws.Range("N102").Value = -4798.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5957625.5
$ws.Range("I132").Value = 7357725.5
$ws.Range("K132").Value = 22073176.5
$ws.Range("M132").Value = -22070646.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4574.8335
$ws.Range("I7").Value = 4299.6665
$ws.Range("J7").Value = 4850
$ws.Range("K7").Value = 4299.6665
$ws.Range("L7").Value = 4850
$ws.Range("M7").Value = -4187.6665
$ws.Range("N7").Value = -5074

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1736.4736
$ws.Range("I16").Value = 1287.3846
$ws.Range("K16").Value = 1287.3846
$ws.Range("M16").Value = -1117.3846

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2680
$ws.Range("I22").Value = 3192.8572
$ws.Range("K22").Value = 3192.8572
$ws.Range("M22").Value = -2897.8572

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2680
$ws.Range("I27").Value = 3192.8572
$ws.Range("K27").Value = 3192.8572
$ws.Range("M27").Value = -3085.8572

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1414.3334
$ws.Range("I46").Value = 1527.6923
$ws.Range("J46").Value = 1119.6
$ws.Range("K46").Value = 1527.6923
$ws.Range("L46").Value = 1119.6
$ws.Range("M46").Value = -1339.6923
$ws.Range("N46").Value = -1495.6

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4609.0557
$ws.Range("I61").Value = 4544.9414
$ws.Range("K61").Value = 4544.9414
$ws.Range("M61").Value = -4342.9414

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4609.0557
$ws.Range("I113").Value = 4544.9414
$ws.Range("K113").Value = 4544.9414
$ws.Range("M113").Value = -2374.9414

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10199.667
$ws.Range("I122").Value = 10199.667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 30599.001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -28149.001
$ws.Range("N122").ClearContents()

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4574.8335
$ws.Range("I126").Value = 4299.6665
$ws.Range("J126").Value = 4850
$ws.Range("K126").Value = 12898.9995
$ws.Range("L126").Value = 14550
$ws.Range("M126").Value = -10428.9995
$ws.Range("N126").Value = -19490

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 761.5
$ws.Range("I107").Value = 732
$ws.Range("K107").Value = 2196
$ws.Range("M107").Value = -276
